$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(9, 8).Value = 79.73684
$ws.Cells.Item(9, 9).Value = 84.375
$ws.Cells.Item(9, 10).Value = 55
$ws.Cells.Item(9, 11).Value = 84.375
$ws.Cells.Item(9, 12).Value = 55
$ws.Cells.Item(9, 13).Value = 84.625
$ws.Cells.Item(9, 14).Value = -393
$ws.Cells.Item(38, 8).Value = 1926.5555
$ws.Cells.Item(38, 10).Value = 8494.5
$ws.Cells.Item(38, 12).Value = 25483.5
$ws.Cells.Item(38, 14).Value = -26227.5
$ws.Cells.Item(75, 8).Value = 0
$ws.Cells.Item(75, 10).Value = 0
$ws.Cells.Item(75, 12).Value = 0
$ws.Cells.Item(75, 14).ClearContents()
$ws.Cells.Item(78, 8).Value = 0
$ws.Cells.Item(78, 10).Value = 0
$ws.Cells.Item(78, 12).Value = 0
$ws.Cells.Item(78, 14).ClearContents()
$ws.Cells.Item(86, 8).Value = 3243.2
$ws.Cells.Item(86, 9).Value = 3245.3333
$ws.Cells.Item(86, 10).Value = 3240
$ws.Cells.Item(86, 11).Value = 3245.3333
$ws.Cells.Item(86, 12).Value = 3240
$ws.Cells.Item(86, 13).Value = -2122.3333
$ws.Cells.Item(86, 14).Value = -5486
$ws.Cells.Item(89, 8).Value = 3243.2
$ws.Cells.Item(89, 9).Value = 3245.3333
$ws.Cells.Item(89, 10).Value = 3240
$ws.Cells.Item(89, 11).Value = 16226.6665
$ws.Cells.Item(89, 12).Value = 16200
$ws.Cells.Item(89, 13).Value = -10610.6665
$ws.Cells.Item(89, 14).Value = -27432
$ws.Cells.Item(107, 8).Value = 499
$ws.Cells.Item(107, 9).Value = 499
$ws.Cells.Item(107, 11).Value = 499
$ws.Cells.Item(107, 13).Value = 1421
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(61, 8).Value = 0
$ws.Cells.Item(61, 9).Value = 0
$ws.Cells.Item(61, 10).Value = 0
$ws.Cells.Item(61, 11).Value = 0
$ws.Cells.Item(61, 12).Value = 0
$ws.Cells.Item(61, 13).ClearContents()
$ws.Cells.Item(61, 14).ClearContents()
$ws.Cells.Item(132, 8).Value = 133.5
$ws.Cells.Item(132, 9).Value = 133.5
$ws.Cells.Item(132, 11).Value = 400.5
$ws.Cells.Item(132, 13).Value = 2129.5
$ws.Cells.Item(136, 8).Value = 0
$ws.Cells.Item(136, 9).Value = 0
$ws.Cells.Item(136, 10).Value = 0
$ws.Cells.Item(136, 11).Value = 0
$ws.Cells.Item(136, 12).Value = 0
$ws.Cells.Item(136, 13).ClearContents()
$ws.Cells.Item(136, 14).ClearContents()
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(29, 8).Value = 4010.3333
$ws.Cells.Item(29, 9).Value = 4010.3333
$ws.Cells.Item(29, 10).Value = 0
$ws.Cells.Item(29, 11).Value = 4010.3333
$ws.Cells.Item(29, 12).Value = 0
$ws.Cells.Item(29, 13).Value = -3721.3333
$ws.Cells.Item(29, 14).ClearContents()
$ws.Cells.Item(57, 8).Value = 100000
$ws.Cells.Item(57, 10).Value = 100000
$ws.Cells.Item(57, 12).Value = 100000
$ws.Cells.Item(57, 14).Value = -101440
$ws.Cells.Item(107, 8).Value = 829.75
$ws.Cells.Item(107, 9).Value = 829.75
$ws.Cells.Item(107, 10).Value = 0
$ws.Cells.Item(107, 11).Value = 829.75
$ws.Cells.Item(107, 12).Value = 0
$ws.Cells.Item(107, 13).Value = 1090.25
$ws.Cells.Item(107, 14).ClearContents()
$ws.Cells.Item(134, 8).Value = 822.3333
$ws.Cells.Item(134, 9).Value = 822.3333
$ws.Cells.Item(134, 10).Value = 0
$ws.Cells.Item(134, 11).Value = 2466.9999
$ws.Cells.Item(134, 12).Value = 0
$ws.Cells.Item(134, 13).Value = 68.0001000000002
$ws.Cells.Item(134, 14).ClearContents()
$ws.Cells.Item(136, 8).Value = 100000
$ws.Cells.Item(136, 10).Value = 100000
$ws.Cells.Item(136, 12).Value = 100000
$ws.Cells.Item(136, 14).Value = -110200
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(12, 8).Value = 0
$ws.Cells.Item(12, 10).Value = 0
$ws.Cells.Item(12, 12).Value = 0
$ws.Cells.Item(12, 14).ClearContents()
$ws.Cells.Item(31, 8).Value = 0
$ws.Cells.Item(31, 10).Value = 0
$ws.Cells.Item(31, 12).Value = 0
$ws.Cells.Item(31, 14).ClearContents()
$ws.Cells.Item(34, 8).Value = 0
$ws.Cells.Item(34, 10).Value = 0
$ws.Cells.Item(34, 12).Value = 0
$ws.Cells.Item(34, 14).ClearContents()
$ws.Cells.Item(35, 8).Value = 6565.2
$ws.Cells.Item(35, 9).Value = 950
$ws.Cells.Item(35, 10).Value = 14988
$ws.Cells.Item(35, 11).Value = 950
$ws.Cells.Item(35, 12).Value = 14988
$ws.Cells.Item(35, 13).Value = -656
$ws.Cells.Item(35, 14).Value = -15576
$ws.Cells.Item(86, 8).Value = 0
$ws.Cells.Item(86, 9).Value = 0
$ws.Cells.Item(86, 10).Value = 0
$ws.Cells.Item(86, 11).Value = 0
$ws.Cells.Item(86, 12).Value = 0
$ws.Cells.Item(86, 13).ClearContents()
$ws.Cells.Item(86, 14).ClearContents()
$ws.Cells.Item(89, 8).Value = 0
$ws.Cells.Item(89, 9).Value = 0
$ws.Cells.Item(89, 10).Value = 0
$ws.Cells.Item(89, 11).Value = 0
$ws.Cells.Item(89, 12).Value = 0
$ws.Cells.Item(89, 13).ClearContents()
$ws.Cells.Item(89, 14).ClearContents()
$ws.Cells.Item(99, 8).Value = 3000
$ws.Cells.Item(99, 9).Value = 3000
$ws.Cells.Item(99, 11).Value = 3000
$ws.Cells.Item(99, 13).Value = -1502
$ws.Cells.Item(126, 8).Value = 3000
$ws.Cells.Item(126, 9).Value = 3000
$ws.Cells.Item(126, 11).Value = 9000
$ws.Cells.Item(126, 13).Value = -6530
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(2, 8).Value = 717.8929000000001
$ws.Cells.Item(2, 9).Value = 359.25
$ws.Cells.Item(2, 10).Value = 2869.75
$ws.Cells.Item(2, 11).Value = 2155.5
$ws.Cells.Item(2, 12).Value = 17218.5
$ws.Cells.Item(2, 13).Value = -2042.5
$ws.Cells.Item(2, 14).Value = -17444.5
$ws.Cells.Item(81, 8).Value = 4853.25
$ws.Cells.Item(81, 10).Value = 700
$ws.Cells.Item(81, 12).Value = 2100
$ws.Cells.Item(81, 14).Value = -4346
$ws.Cells.Item(84, 8).Value = 4853.25
$ws.Cells.Item(84, 10).Value = 700
$ws.Cells.Item(84, 12).Value = 6300
$ws.Cells.Item(84, 14).Value = -17532
$ws.Cells.Item(92, 8).Value = 0
$ws.Cells.Item(92, 9).Value = 0
$ws.Cells.Item(92, 10).Value = 0
$ws.Cells.Item(92, 11).Value = 0
$ws.Cells.Item(92, 12).Value = 0
$ws.Cells.Item(92, 13).ClearContents()
$ws.Cells.Item(92, 14).ClearContents()
$ws.Cells.Item(121, 8).Value = 382.5
$ws.Cells.Item(121, 10).Value = 50
$ws.Cells.Item(121, 12).Value = 150
$ws.Cells.Item(121, 14).Value = -2770
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(52, 8).Value = 40000
$ws.Cells.Item(52, 9).Value = 0
$ws.Cells.Item(52, 11).Value = 0
$ws.Cells.Item(52, 13).ClearContents()
$ws.Cells.Item(97, 8).Value = 734.2857
$ws.Cells.Item(97, 9).Value = 746.6667
$ws.Cells.Item(97, 10).Value = 725
$ws.Cells.Item(97, 11).Value = 746.6667
$ws.Cells.Item(97, 12).Value = 725
$ws.Cells.Item(97, 13).Value = -250.6667
$ws.Cells.Item(97, 14).Value = -1717
$ws.Cells.Item(113, 8).Value = 2288.5
$ws.Cells.Item(113, 9).Value = 2288.5
$ws.Cells.Item(113, 11).Value = 2288.5
$ws.Cells.Item(113, 13).Value = -118.5
$ws.Cells.Item(132, 8).Value = 3970.3333
$ws.Cells.Item(132, 9).Value = 4205.5
$ws.Cells.Item(132, 10).Value = 3500
$ws.Cells.Item(132, 11).Value = 12616.5
$ws.Cells.Item(132, 12).Value = 10500
$ws.Cells.Item(132, 13).Value = -10086.5
$ws.Cells.Item(132, 14).Value = -15560
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(46, 8).Value = 4459.4
$ws.Cells.Item(46, 9).Value = 5600
$ws.Cells.Item(46, 10).Value = 2748.5
$ws.Cells.Item(46, 11).Value = 5600
$ws.Cells.Item(46, 12).Value = 2748.5
$ws.Cells.Item(46, 13).Value = -5412
$ws.Cells.Item(46, 14).Value = -3124.5
$ws.Cells.Item(56, 8).Value = 16499
$ws.Cells.Item(56, 9).Value = 8000
$ws.Cells.Item(56, 11).Value = 8000
$ws.Cells.Item(56, 13).Value = -7309
$ws.Cells.Item(82, 8).Value = 2329.6667
$ws.Cells.Item(82, 9).Value = 994.5
$ws.Cells.Item(82, 11).Value = 994.5
$ws.Cells.Item(82, 13).Value = -633.5
$ws.Cells.Item(85, 8).Value = 2329.6667
$ws.Cells.Item(85, 9).Value = 994.5
$ws.Cells.Item(85, 11).Value = 994.5
$ws.Cells.Item(85, 13).Value = 253.5
$ws.Cells.Item(95, 8).Value = 50344
$ws.Cells.Item(95, 10).Value = 50344
$ws.Cells.Item(95, 12).Value = 50344
$ws.Cells.Item(95, 14).Value = -55836
$ws.Cells.Item(104, 8).Value = 0
$ws.Cells.Item(104, 10).Value = 0
$ws.Cells.Item(104, 12).Value = 0
$ws.Cells.Item(104, 14).ClearContents()
$ws.Cells.Item(122, 8).Value = 1199
$ws.Cells.Item(122, 9).Value = 1199
$ws.Cells.Item(122, 11).Value = 3597
$ws.Cells.Item(122, 13).Value = -1147
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(19, 8).Value = 65000000
$ws.Cells.Item(19, 9).Value = 100000000
$ws.Cells.Item(19, 10).Value = 30000000
$ws.Cells.Item(19, 11).Value = 100000000
$ws.Cells.Item(19, 12).Value = 30000000
$ws.Cells.Item(19, 13).Value = -99999826
$ws.Cells.Item(19, 14).Value = -30000348
$ws.Cells.Item(68, 8).Value = 0
$ws.Cells.Item(68, 10).Value = 0
$ws.Cells.Item(68, 12).Value = 0
$ws.Cells.Item(68, 14).ClearContents()
$ws.Cells.Item(71, 8).Value = 0
$ws.Cells.Item(71, 10).Value = 0
$ws.Cells.Item(71, 12).Value = 0
$ws.Cells.Item(71, 14).ClearContents()
$ws.Cells.Item(82, 8).Value = 44999.5
$ws.Cells.Item(82, 10).Value = 44999.5
$ws.Cells.Item(82, 12).Value = 44999.5
$ws.Cells.Item(82, 14).Value = -45765.5
$ws.Cells.Item(85, 8).Value = 44999.5
$ws.Cells.Item(85, 10).Value = 44999.5
$ws.Cells.Item(85, 12).Value = 44999.5
$ws.Cells.Item(85, 14).Value = -47651.5
$ws.Cells.Item(98, 8).Value = 39600
$ws.Cells.Item(98, 10).Value = 39600
$ws.Cells.Item(98, 12).Value = 39600
$ws.Cells.Item(98, 14).Value = -45590
$ws.Cells.Item(104, 8).Value = 22821.666
$ws.Cells.Item(104, 10).Value = 22821.666
$ws.Cells.Item(104, 12).Value = 22821.666
$ws.Cells.Item(104, 14).Value = -29809.666
$ws.Cells.Item(126, 8).Value = 304
$ws.Cells.Item(126, 9).Value = 304
$ws.Cells.Item(126, 10).Value = 0
$ws.Cells.Item(126, 11).Value = 912
$ws.Cells.Item(126, 12).Value = 0
$ws.Cells.Item(126, 13).Value = 1558
$ws.Cells.Item(126, 14).ClearContents()
$ws.Cells.Item(132, 8).Value = 4393.6665
$ws.Cells.Item(132, 9).Value = 4393.6665
$ws.Cells.Item(132, 11).Value = 13180.9995
$ws.Cells.Item(132, 13).Value = -10650.9995
